# Auto-generated edit script
# Applies numeric cell-value updates to match the target diff.
# Cells whose new value is absent in the diff are cleared (ClearContents);
# cells that are new in the diff are written for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1022.7692
$ws.Range("I28").Value = 1091.3334
$ws.Range("K28").Value = 1091.3334
$ws.Range("M28").Value = -606.3334
$ws.Range("H33").Value = 268.0909
$ws.Range("I33").Value = 172.11111
$ws.Range("K33").Value = 172.11111
$ws.Range("M33").Value = 56.88889
$ws.Range("H70").Value = 9624.25
$ws.Range("I70").Value = 11199.4
$ws.Range("J70").Value = 6999
$ws.Range("K70").Value = 33598.2
$ws.Range("L70").Value = 20997
$ws.Range("M70").Value = -33328.2
$ws.Range("N70").Value = -21537
$ws.Range("H73").Value = 9624.25
$ws.Range("I73").Value = 11199.4
$ws.Range("J73").Value = 6999
$ws.Range("K73").Value = 33598.2
$ws.Range("L73").Value = 20997
$ws.Range("M73").Value = -32662.2
$ws.Range("N73").Value = -22869
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H80").Value = 1070.6666
$ws.Range("J80").Value = 1697
$ws.Range("L80").Value = 5091
$ws.Range("N80").Value = -7087
$ws.Range("H83").Value = 1070.6666
$ws.Range("J83").Value = 1697
$ws.Range("L83").Value = 15273
$ws.Range("N83").Value = -25257
$ws.Range("H112").Value = 2462.0715
$ws.Range("I112").Value = 1241
$ws.Range("K112").Value = 3723
$ws.Range("M112").Value = -2615
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11008
$ws.Range("H115").Value = 1467
$ws.Range("I115").Value = 635
$ws.Range("J115").Value = 2299
$ws.Range("K115").Value = 1905
$ws.Range("L115").Value = 6897
$ws.Range("M115").Value = -338
$ws.Range("N115").Value = -10031
$ws.Range("H119").Value = 1400
$ws.Range("J119").Value = 1400
$ws.Range("L119").Value = 4200
$ws.Range("N119").Value = -13876
$ws.Range("H132").Value = 2417
$ws.Range("I132").Value = 2218.3076
$ws.Range("K132").Value = 6654.9228
$ws.Range("M132").Value = -4124.9228
$ws.Range("H141").Value = 2135.1
$ws.Range("I141").Value = 1457.125
$ws.Range("J141").Value = 4847
$ws.Range("K141").Value = 4371.375
$ws.Range("L141").Value = 14541
$ws.Range("M141").Value = 808.625
$ws.Range("N141").Value = -24901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2753934.5
$ws.Range("I32").Value = 2596672.8
$ws.Range("K32").Value = 2596672.8
$ws.Range("M32").Value = -2596385.8
$ws.Range("H97").Value = 1175.7142
$ws.Range("I97").Value = 1176.909
$ws.Range("K97").Value = 1176.909
$ws.Range("M97").Value = -680.9090000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 73
$ws.Range("I22").Value = 73
$ws.Range("J22").Value = 73
$ws.Range("K22").Value = 73
$ws.Range("L22").Value = 73
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -419
$ws.Range("H86").Value = 6273.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 6273.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 2664.75
$ws.Range("I99").Value = 2563.8
$ws.Range("K99").Value = 2563.8
$ws.Range("M99").Value = -1065.8
$ws.Range("H107").Value = 2063.4546
$ws.Range("I107").Value = 1833
$ws.Range("K107").Value = 1833
$ws.Range("M107").Value = 87

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1034.2
$ws.Range("I31").Value = 852.3333
$ws.Range("J31").Value = 1136.5
$ws.Range("K31").Value = 852.3333
$ws.Range("L31").Value = 1136.5
$ws.Range("M31").Value = -557.3333
$ws.Range("N31").Value = -1726.5
$ws.Range("H34").Value = 1034.2
$ws.Range("I34").Value = 852.3333
$ws.Range("J34").Value = 1136.5
$ws.Range("K34").Value = 852.3333
$ws.Range("L34").Value = 1136.5
$ws.Range("M34").Value = -650.3333
$ws.Range("N34").Value = -1540.5
$ws.Range("H105").Value = 3213
$ws.Range("I105").Value = 2516.1667
$ws.Range("J105").Value = 4049.2
$ws.Range("K105").Value = 2516.1667
$ws.Range("L105").Value = 4049.2
$ws.Range("M105").Value = -769.1667000000002
$ws.Range("N105").Value = -7543.2
$ws.Range("H107").Value = 852.5
$ws.Range("I107").Value = 748.25
$ws.Range("K107").Value = 748.25
$ws.Range("M107").Value = 1171.75
$ws.Range("H122").Value = 1200.2778
$ws.Range("I122").Value = 1401.5555
$ws.Range("K122").Value = 4204.666499999999
$ws.Range("M122").Value = -1754.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5694824
$ws.Range("I4").Value = 6114527.5
$ws.Range("J4").Value = 3176603.2
$ws.Range("K4").Value = 18343582.5
$ws.Range("L4").Value = 9529809.600000001
$ws.Range("M4").Value = -18343470.5
$ws.Range("N4").Value = -9530033.600000001
$ws.Range("H11").Value = 246122.58
$ws.Range("I11").Value = 246122.58
$ws.Range("K11").Value = 738367.74
$ws.Range("M11").Value = -738227.74
$ws.Range("H34").Value = 1000
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 3000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -3354
$ws.Range("H132").Value = 1842.2858
$ws.Range("I132").Value = 1982.4
$ws.Range("K132").Value = 17841.6
$ws.Range("M132").Value = -15311.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 237500
$ws.Range("I5").Value = 225000
$ws.Range("K5").Value = 225000
$ws.Range("M5").Value = -224888
$ws.Range("H7").Value = 750000
$ws.Range("I7").Value = 1000000
$ws.Range("J7").Value = 500000
$ws.Range("K7").Value = 1000000
$ws.Range("L7").Value = 500000
$ws.Range("M7").Value = -999888
$ws.Range("N7").Value = -500224
$ws.Range("H8").Value = 750000
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 500000
$ws.Range("K8").Value = 1000000
$ws.Range("L8").Value = 500000
$ws.Range("M8").Value = -999861
$ws.Range("N8").Value = -500278
$ws.Range("H102").Value = 2299
$ws.Range("I102").Value = 1984.8572
$ws.Range("J102").Value = 2848.75
$ws.Range("K102").Value = 1984.8572
$ws.Range("L102").Value = 2848.75
$ws.Range("M102").Value = -362.8571999999999
$ws.Range("N102").Value = -6092.75
$ws.Range("H126").Value = 5496
$ws.Range("I126").Value = 4231
$ws.Range("J126").Value = 5812.25
$ws.Range("K126").Value = 12693
$ws.Range("L126").Value = 17436.75
$ws.Range("M126").Value = -10223
$ws.Range("N126").Value = -22376.75
$ws.Range("H132").Value = 4001
$ws.Range("I132").Value = 4001
$ws.Range("K132").Value = 12003
$ws.Range("M132").Value = -9473

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 40006
$ws.Range("I23").Value = 40006
$ws.Range("K23").Value = 40006
$ws.Range("M23").Value = -39776
$ws.Range("H40").Value = 3003.3333
$ws.Range("I40").Value = 2036
$ws.Range("K40").Value = 2036
$ws.Range("M40").Value = -1900
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 3388.3
$ws.Range("I132").Value = 3246.7693
$ws.Range("J132").Value = 3651.1428
$ws.Range("K132").Value = 9740.3079
$ws.Range("L132").Value = 10953.4284
$ws.Range("M132").Value = -7210.3079
$ws.Range("N132").Value = -16013.4284
$ws.Range("H136").Value = 1413.1
$ws.Range("I136").Value = 1350
$ws.Range("K136").Value = 4050
$ws.Range("M136").Value = -1500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 880000
$ws.Range("I2").Value = 10000
$ws.Range("K2").Value = 10000
$ws.Range("M2").Value = -9888
$ws.Range("H4").Value = 10018194
$ws.Range("I4").Value = 25000300
$ws.Range("J4").Value = 30123.334
$ws.Range("K4").Value = 25000300
$ws.Range("L4").Value = 30123.334
$ws.Range("M4").Value = -25000187
$ws.Range("N4").Value = -30349.334
$ws.Range("H96").Value = 4409.2
$ws.Range("I96").Value = 4349
$ws.Range("K96").Value = 4349
$ws.Range("M96").Value = -2976
$ws.Range("H122").Value = 3781.682
$ws.Range("J122").Value = 5059.125
$ws.Range("L122").Value = 15177.375
$ws.Range("N122").Value = -20077.375
$ws.Range("H132").Value = 2068.9
$ws.Range("I132").Value = 1898.625
$ws.Range("K132").Value = 5695.875
$ws.Range("M132").Value = -3165.875
